$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (new TC data)
$ws.Range("A2").Value = "hcszK894"
$ws.Range("B2").Value = 231011229
$ws.Range("C2").Value = "lovvxxj63"
$ws.Range("D2").Value = "m`$2W#cU9"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "uQwOJFdm"
$ws.Range("G2").Value = "NxLV"
$ws.Range("H2").Value = "Candidate"

# Add new row 3, matching the bordered style used by row 2
$ws.Range("A3:H3").Borders.LineStyle = 1
$ws.Range("A3:H3").Borders.Weight = 2

$ws.Range("A3").Value = "zjEbF698"
$ws.Range("B3").Value = 231011228
$ws.Range("C3").Value = "fnnjlff28"
$ws.Range("D3").Value = "xJ65!W#y"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "jnkrqvHy"
$ws.Range("G3").Value = "GkoO"
$ws.Range("H3").Value = "Candidate"
